# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# OFF sheet - Row 3 ("R") gets updated with Wild Card round stats added in
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 450
$wsOff.Range("C3").Value = 333
$wsOff.Range("D3").Value = 112
$wsOff.Range("E3").Value = 52

# DEF sheet - Row 3 ("R") gets updated with Wild Card round stats added in
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 474
$wsDef.Range("C3").Value = 314
$wsDef.Range("D3").Value = 113
$wsDef.Range("E3").Value = 58
